# Refresh the crypto price/volume snapshot (scheduled GitHub Actions update).
# Maps each changed cell reference to its new literal text value.
$updates = [ordered]@{
    'D2' = '68.124.28'
    'D3' = '3.798.28'
    'E3' = '  -0.22%  '
    'D4' = '0.999'
    'E4' = '  -0.20%  '
    'E5' = '  +0.70%  '
    'D6' = '165.22'
    'E6' = '  -1.51%  '
    'E7' = '  -0.04%  '
    'E8' = '  -0.55%  '
    'E9' = '  -0.90%  '
    'E10' = '  +0.43%  '
    'D11' = '6.48'
    'E11' = '  +3.01%  '
    'E12' = '  -2.02%  '
    'D13' = '35.91'
    'E13' = '  -0.31%  '
    'D14' = '4.431.52'
    'E14' = '  -0.36%  '
    'D15' = '3.766.32'
    'E15' = '  -0.63%  '
    'D16' = '68.117.19'
    'E17' = '  -0.96%  '
    'E18' = '  +2.36%  '
    'D19' = '7.09'
    'E19' = '  -0.13%  '
    'D20' = '461.64'
    'E20' = '  -0.05%  '
    'D21' = '9.72'
    'E21' = '  -2.36%  '
    'E22' = '  +0.07%  '
    'E23' = '  -4.55%  '
    'D24' = '83.02'
    'E24' = '  -0.70%  '
    'D25' = '12.03'
    'E25' = '  -0.53%  '
    'E26' = '  +0.18%  '
    'B27' = 'RenderToken'
    'C27' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D27' = '10.03'
    'E27' = '  +0.09%  '
    'B28' = 'Dai'
    'C28' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D28' = '1.00'
    'E28' = '  +0.08%  '
    'D29' = '3.943.72'
    'E29' = '  -0.32%  '
    'B30' = 'ImmutableX'
    'C30' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D30' = '2.24'
    'E30' = '  -0.30%  '
    'B31' = 'PancakeSwap'
    'C31' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D31' = '2.65'
    'E31' = '  -4.57%  '
    'E32' = '  +0.70%  '
    'D33' = '29.36'
    'E33' = '  -1.07%  '
    'D34' = '0.999'
    'E34' = '  +0.15%  '
    'D35' = '9.05'
    'E35' = '  -0.47%  '
    'D36' = '0.0998'
    'E36' = '  -0.14%  '
    'D37' = '3.35'
    'E37' = '  -2.26%  '
    'E38' = '  +0.92%  '
    'E39' = '  +1.00%  '
    'D40' = '0.989'
    'E40' = '  -1.42%  '
    'E41' = '  -0.02%  '
    'E42' = '  +0.00%  '
    'D43' = '47.63'
    'E43' = '  -1.07%  '
    'E44' = '  +0.05%  '
    'E45' = '  -1.22%  '
    'D46' = '152.56'
    'E46' = '  +2.88%  '
    'D47' = '8.37'
    'E47' = '  +0.50%  '
    'D48' = '1.88'
    'E48' = '  +1.87%  '
    'E49' = '  +1.23%  '
    'D50' = '392.35'
    'E50' = '  -0.96%  '
    'D51' = '26.56'
    'E51' = '  -1.64%  '
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)

    # Columns D/E hold text-formatted numbers/percentages (e.g. "68.124.28",
    # "0.999", "  -0.22%  "). Force text storage for anything that would
    # otherwise be auto-parsed as a plain number so it round-trips as a
    # literal string, matching the source feed's formatting.
    if ($ref.StartsWith('D') -and ($value -match '^-?\d+(\.\d+)?$')) {
        $cell.NumberFormat = "@"
    }

    $cell.Value = $value
}
